$wb = $excel.ActiveWorkbook

# --- Sheet "ev_charging_uc": update the two lookup-table strings (D11/D12) ---
$wsUc = $wb.Worksheets.Item("ev_charging_uc")
$wsUc.Range("D11").Value2 = "S1aH2,S2aH3,S1aH3,S2aH2"
$wsUc.Range("D12").Value2 = "S2aH1,S1aH1,S2aH4,S1aH4"

# --- Sheet "ts12_clu": update commodity labels, season labels and swapped values ---
$wsTs = $wb.Worksheets.Item("ts12_clu")

# AG11:AG18 commodity label changed from "ELC" to "Elec"
$wsTs.Range("AG11").Value2 = "Elec"
$wsTs.Range("AG12").Value2 = "Elec"
$wsTs.Range("AG13").Value2 = "Elec"
$wsTs.Range("AG14").Value2 = "Elec"
$wsTs.Range("AG15").Value2 = "Elec"
$wsTs.Range("AG16").Value2 = "Elec"
$wsTs.Range("AG17").Value2 = "Elec"
$wsTs.Range("AG18").Value2 = "Elec"

# AK11/AK12 season labels swapped, AL11/AL12 values swapped
$wsTs.Range("AK11").Value2 = "S1"
$wsTs.Range("AK12").Value2 = "S2"
$wsTs.Range("AL11").Value2 = 1.0373322535863025
$wsTs.Range("AL12").Value2 = 0.16266774641369736

# Recalculate so dependent formulas (e.g. HLOOKUP results) refresh their cached values
$wb.Application.Calculate()
